$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 131046824
$ws.Range("Q2").Value = 401653
$ws.Range("R2").Value = 6818054
$ws.Range("Z2").Value = "14:50"
$ws.Range("AB2").Value = "14:50"

# --- Row 3 ---
$ws.Range("A3").Value = 131046825
$ws.Range("Q3").Value = 401650
$ws.Range("R3").Value = 6818017
$ws.Range("Z3").Value = "14:52"
$ws.Range("AB3").Value = "14:52"

# --- Row 8 ---
$ws.Range("A8").Value = 131046822
$ws.Range("B8").Value = 79243
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("M8").ClearContents()
$ws.Range("Q8").Value = 401699
$ws.Range("R8").Value = 6818070
$ws.Range("Z8").Value = "14:49"
$ws.Range("AB8").Value = "14:49"
$ws.Range("AC8").ClearContents()
$ws.Range("AE8").Value = $false

# --- Row 9 ---
$ws.Range("A9").Value = 131047013
$ws.Range("B9").Value = 57884
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("M9").Value = "färska spår"
$ws.Range("Q9").Value = 401631
$ws.Range("R9").Value = 6817903
$ws.Range("Z9").Value = "14:57"
$ws.Range("AB9").Value = "14:57"
$ws.Range("AC9").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE9").Value = $true

# --- Row 10 ---
$ws.Range("A10").Value = 131046823
$ws.Range("B10").Value = 79243
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("M10").ClearContents()
$ws.Range("Q10").Value = 401661
$ws.Range("R10").Value = 6818064
$ws.Range("Z10").Value = "14:50"
$ws.Range("AB10").Value = "14:50"
$ws.Range("AC10").ClearContents()

# --- Row 11 ---
$ws.Range("A11").Value = 131046773
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("M11").Value = "äldre spår"
$ws.Range("Q11").Value = 401346
$ws.Range("R11").Value = 6818162
$ws.Range("Z11").Value = "15:23"
$ws.Range("AB11").Value = "15:23"
$ws.Range("AC11").Value = "Äldre ringhack (gran)"

# --- Row 23 ---
$ws.Range("A23").Value = 131046833
$ws.Range("Q23").Value = 401322
$ws.Range("R23").Value = 6818367
$ws.Range("Z23").Value = "15:29"
$ws.Range("AB23").Value = "15:29"

# --- Row 24 ---
$ws.Range("A24").Value = 131046831
$ws.Range("Q24").Value = 401378
$ws.Range("R24").Value = 6818089
$ws.Range("Z24").Value = "15:21"
$ws.Range("AB24").Value = "15:21"

# --- Row 25 ---
$ws.Range("A25").Value = 131047014
$ws.Range("B25").Value = 57884
$ws.Range("E25").Value = 100109
$ws.Range("F25").Value = "Tretåig hackspett"
$ws.Range("G25").Value = "Picoides tridactylus"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("M25").Value = "färska spår"
$ws.Range("Q25").Value = 401378
$ws.Range("R25").Value = 6818082
$ws.Range("Z25").Value = "15:21"
$ws.Range("AB25").Value = "15:21"
$ws.Range("AC25").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE25").Value = $true

# --- Row 26 ---
$ws.Range("A26").Value = 131046832
$ws.Range("Q26").Value = 401350
$ws.Range("R26").Value = 6818162
$ws.Range("Z26").Value = "15:24"
$ws.Range("AB26").Value = "15:24"

# --- Row 27 ---
$ws.Range("A27").Value = 131046826
$ws.Range("B27").Value = 79243
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("M27").ClearContents()
$ws.Range("Q27").Value = 401647
$ws.Range("R27").Value = 6817965
$ws.Range("Z27").Value = "14:54"
$ws.Range("AB27").Value = "14:54"
$ws.Range("AC27").ClearContents()
$ws.Range("AE27").Value = $false

